# ND02.xlsx edit: "Merge back T2A sheet in the test files"
# Adds a new worksheet named "T2A" at the end of the workbook, containing
# the same 41-column layout as WMT_Extract, populated with two data rows.

$wb = $excel.ActiveWorkbook

# 1. Add the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "T2A"

# 2. Apply formatting that mirrors existing styles already used elsewhere in
#    the workbook, so no redundant styles/fonts get created.

# Header row (row 1) -> same style as WMT_Extract!Y1 (font "Arial" 13pt black)
$wmtExtract = $wb.Worksheets.Item("WMT_Extract")
$wmtExtract.Range("Y1").Copy()
$ws.Range("A1:AO1").PasteSpecial(-4122)

# Body rows (rows 2-3), columns A:AN -> same style as ARMS!D2 (font "Calibri" 12pt black)
$arms = $wb.Worksheets.Item("ARMS")
$arms.Range("D2").Copy()
$ws.Range("A2:AN3").PasteSpecial(-4122)

# Body rows (rows 2-3), column AO (Datestamp) -> same base style, with a date+time number format
$arms.Range("D2").Copy()
$ws.Range("AO2:AO3").PasteSpecial(-4122)
$ws.Range("AO2:AO3").NumberFormat = "m/d/yy h:mm"

$excel.CutCopyMode = 0

# 3. Populate the cell values.
$data = New-Object 'object[,]' 3,41
$data[0,0] = "Trust"
$data[0,1] = "Region_Desc"
$data[0,2] = "Region_Code"
$data[0,3] = "Ldu_Desc"
$data[0,4] = "Ldu_Code"
$data[0,5] = "Team_Desc"
$data[0,6] = "Team_Code"
$data[0,7] = "OM_Surname"
$data[0,8] = "OM_Forename"
$data[0,9] = "OM_Grade_Code"
$data[0,10] = "OM_Key"
$data[0,11] = "CommTier0"
$data[0,12] = "CommTierD2"
$data[0,13] = "CommTier1a"
$data[0,14] = "CommTier"
$data[0,15] = "CommTierD1"
$data[0,16] = "CommTierC2"
$data[0,17] = "CommTier3a"
$data[0,18] = "CommTierC1"
$data[0,19] = "CommTierB2"
$data[0,20] = "CommTierB1"
$data[0,21] = "CommTierA"
$data[0,22] = "LicenceTier0"
$data[0,23] = "LicenceTierD2"
$data[0,24] = "LicenceTierD1"
$data[0,25] = "LicenceTierC2"
$data[0,26] = "LicenceTierC1"
$data[0,27] = "LicenceTierB2"
$data[0,28] = "LicenceTierB1"
$data[0,29] = "LicenceTierA"
$data[0,30] = "CustTier0"
$data[0,31] = "CustTierD2"
$data[0,32] = "CustTierD1"
$data[0,33] = "CustTierC2"
$data[0,34] = "CustTierC1"
$data[0,35] = "CustTierB2"
$data[0,36] = "CustTierB1"
$data[0,37] = "CustTierA"
$data[0,38] = "ComIn1st16Weeks"
$data[0,39] = "LicIn1st16Weeks"
$data[0,40] = "Datestamp"
$data[1,0] = "Farringdon"
$data[1,1] = "London II"
$data[1,2] = "ND02"
$data[1,3] = "KainosLDU"
$data[1,4] = "KNS"
$data[1,5] = "WMT Team"
$data[1,6] = "WMT (ND02)"
$data[1,7] = "Swann II"
$data[1,8] = "Tom"
$data[1,9] = "C"
$data[1,10] = 1004
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 0
$data[1,14] = 0
$data[1,15] = 0
$data[1,16] = 0
$data[1,17] = 0
$data[1,18] = 0
$data[1,19] = 0
$data[1,20] = 0
$data[1,21] = 0
$data[1,22] = 6
$data[1,23] = 0
$data[1,24] = 0
$data[1,25] = 0
$data[1,26] = 0
$data[1,27] = 0
$data[1,28] = 0
$data[1,29] = 0
$data[1,30] = 0
$data[1,31] = 0
$data[1,32] = 0
$data[1,33] = 0
$data[1,34] = 0
$data[1,35] = 0
$data[1,36] = 0
$data[1,37] = 0
$data[1,38] = 0
$data[1,39] = 0
$data[1,40] = 42795.628472222219
$data[2,0] = "Farringdon"
$data[2,1] = "London II"
$data[2,2] = "ND02"
$data[2,3] = "KainosLDU"
$data[2,4] = "KNS"
$data[2,5] = "WMT Team"
$data[2,6] = "WMT (ND02)"
$data[2,7] = "Wright II"
$data[2,8] = "Andy"
$data[2,9] = "Z"
$data[2,10] = 1005
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 0
$data[2,14] = 0
$data[2,15] = 0
$data[2,16] = 0
$data[2,17] = 0
$data[2,18] = 0
$data[2,19] = 12
$data[2,20] = 0
$data[2,21] = 0
$data[2,22] = 0
$data[2,23] = 0
$data[2,24] = 0
$data[2,25] = 0
$data[2,26] = 0
$data[2,27] = 0
$data[2,28] = 0
$data[2,29] = 0
$data[2,30] = 0
$data[2,31] = 0
$data[2,32] = 0
$data[2,33] = 0
$data[2,34] = 0
$data[2,35] = 0
$data[2,36] = 0
$data[2,37] = 0
$data[2,38] = 0
$data[2,39] = 0
$data[2,40] = 42795.628472222219

$ws.Range("A1:AO3").Value = $data

# 4. Select the full populated range and make this the active sheet/tab,
#    matching the saved view state of the source workbook.
$ws.Range("A1:AO3").Select()
$ws.Activate()
